$d = $word.ActiveDocument

$d.Content.Find.Execute("945×4=3780", $true, $false, $false, $false, $false, $true, 1, $false, "794×5=3970", 2) | Out-Null
$d.Content.Find.Execute("541×4=2164", $true, $false, $false, $false, $false, $true, 1, $false, "419×7=2933", 2) | Out-Null
$d.Content.Find.Execute("263×8=2104", $true, $false, $false, $false, $false, $true, 1, $false, "634×3=1902", 2) | Out-Null
$d.Content.Find.Execute("958×6=5748", $true, $false, $false, $false, $false, $true, 1, $false, "426×6=2556", 2) | Out-Null
$d.Content.Find.Execute("251×9=2259", $true, $false, $false, $false, $false, $true, 1, $false, "141×7=987", 2) | Out-Null
$d.Content.Find.Execute("543×8=4344", $true, $false, $false, $false, $false, $true, 1, $false, "596×8=4768", 2) | Out-Null
$d.Content.Find.Execute("279×3=837", $true, $false, $false, $false, $false, $true, 1, $false, "842×3=2526", 2) | Out-Null
$d.Content.Find.Execute("217×6=1302", $true, $false, $false, $false, $false, $true, 1, $false, "763×4=3052", 2) | Out-Null
$d.Content.Find.Execute("498×3=1494", $true, $false, $false, $false, $false, $true, 1, $false, "635×6=3810", 2) | Out-Null
$d.Content.Find.Execute("693×3=2079", $true, $false, $false, $false, $false, $true, 1, $false, "661×7=4627", 2) | Out-Null
$d.Content.Find.Execute("249×7=1743", $true, $false, $false, $false, $false, $true, 1, $false, "569×3=1707", 2) | Out-Null
$d.Content.Find.Execute("426×5=2130", $true, $false, $false, $false, $false, $true, 1, $false, "591×7=4137", 2) | Out-Null
$d.Content.Find.Execute("990×6=5940", $true, $false, $false, $false, $false, $true, 1, $false, "748×6=4488", 2) | Out-Null
$d.Content.Find.Execute("873×5=4365", $true, $false, $false, $false, $false, $true, 1, $false, "883×5=4415", 2) | Out-Null
$d.Content.Find.Execute("866×4=3464", $true, $false, $false, $false, $false, $true, 1, $false, "422×4=1688", 2) | Out-Null
$d.Content.Find.Execute("577×6=3462", $true, $false, $false, $false, $false, $true, 1, $false, "749×7=5243", 2) | Out-Null
$d.Content.Find.Execute("746×8=5968", $true, $false, $false, $false, $false, $true, 1, $false, "614×4=2456", 2) | Out-Null
$d.Content.Find.Execute("249×2=498", $true, $false, $false, $false, $false, $true, 1, $false, "928×3=2784", 2) | Out-Null
$d.Content.Find.Execute("601×3=1803", $true, $false, $false, $false, $false, $true, 1, $false, "720×8=5760", 2) | Out-Null
$d.Content.Find.Execute("570×6=3420", $true, $false, $false, $false, $false, $true, 1, $false, "508×6=3048", 2) | Out-Null
$d.Content.Find.Execute("169×6=1014", $true, $false, $false, $false, $false, $true, 1, $false, "787×2=1574", 2) | Out-Null
$d.Content.Find.Execute("287×7=2009", $true, $false, $false, $false, $false, $true, 1, $false, "329×3=987", 2) | Out-Null
$d.Content.Find.Execute("991×6=5946", $true, $false, $false, $false, $false, $true, 1, $false, "809×2=1618", 2) | Out-Null
$d.Content.Find.Execute("778×2=1556", $true, $false, $false, $false, $false, $true, 1, $false, "151×4=604", 2) | Out-Null
$d.Content.Find.Execute("955×2=1910", $true, $false, $false, $false, $false, $true, 1, $false, "873×5=4365", 2) | Out-Null
